$d = $word.ActiveDocument

# 1) "N° 00{{n_anuncio}}-2025-MDP-GLDE" -> "N° 00{{n_anuncio}}-2026-MDP/GLDE-SGLCA"
$d.Content.Find.Execute("-2025-MDP-GLDE", $false, $false, $false, $false, $false, $true, 1, $false, "-2026-MDP/GLDE-SGLCA", 2)

# 2) "EXPEDIENTE: {{num_ds}}-2025 (...)" -> "EXPEDIENTE: {{num_ds}}-2026 (...)"
$d.Content.Find.Execute("{{num_ds}}-2025", $false, $false, $false, $false, $false, $true, 1, $false, "{{num_ds}}-2026", 2)
